# Regex de padrão de endereço
# Remove the trailing "Inscrição"/"Matrícula" code (duplicated from column B)
# plus its trailing space that was appended at the end of each address
# string in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $addr = $ws.Cells.Item($r, 3).Value()
    $code = $ws.Cells.Item($r, 2).Value()

    if ($null -ne $addr -and $null -ne $code) {
        $suffix = [string]$code + " "
        if ($addr.EndsWith($suffix)) {
            $newAddr = $addr.Substring(0, $addr.Length - $suffix.Length)
            $ws.Cells.Item($r, 3).Value = $newAddr
        }
    }
}
